$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 125067.5
$ws.Range("I4").Value = 250040
$ws.Range("K4").Value = 250040
$ws.Range("M4").Value = -249926
$ws.Range("H5").Value = 83.333336
$ws.Range("I5").Value = 59.6
$ws.Range("J5").Value = 202
$ws.Range("K5").Value = 59.6
$ws.Range("L5").Value = 202
$ws.Range("M5").Value = 55.4
$ws.Range("N5").Value = -432
$ws.Range("H8").Value = 50
$ws.Range("I8").Value = 50
$ws.Range("K8").Value = 150
$ws.Range("M8").Value = -11
$ws.Range("H9").Value = 137
$ws.Range("I9").Value = 103.083336
$ws.Range("J9").Value = 218.4
$ws.Range("K9").Value = 103.083336
$ws.Range("L9").Value = 218.4
$ws.Range("M9").Value = 65.916664
$ws.Range("N9").Value = -556.4
$ws.Range("H12").Value = 91.888885
$ws.Range("I12").Value = 91.333336
$ws.Range("J12").Value = 93
$ws.Range("K12").Value = 91.333336
$ws.Range("L12").Value = 93
$ws.Range("M12").Value = 78.666664
$ws.Range("N12").Value = -433
$ws.Range("H51").Value = 7620.96
$ws.Range("I51").Value = 1401.5625
$ws.Range("J51").Value = 18677.666
$ws.Range("K51").Value = 1401.5625
$ws.Range("L51").Value = 18677.666
$ws.Range("M51").Value = -917.5625
$ws.Range("N51").Value = -19645.666
$ws.Range("H74").Value = 3373.4
$ws.Range("I74").Value = 3000
$ws.Range("J74").Value = 3622.3333
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 3622.3333
$ws.Range("M74").Value = -2064
$ws.Range("N74").Value = -5494.3333
$ws.Range("H77").Value = 3373.4
$ws.Range("I77").Value = 3000
$ws.Range("J77").Value = 3622.3333
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 18111.6665
$ws.Range("M77").Value = -10320
$ws.Range("N77").Value = -27471.6665
$ws.Range("H132").Value = 6174602.5
$ws.Range("I132").Value = 1848.6
$ws.Range("J132").Value = 37038372
$ws.Range("K132").Value = 5545.799999999999
$ws.Range("L132").Value = 111115116
$ws.Range("M132").Value = -3015.799999999999
$ws.Range("N132").Value = -111120176
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5025761
$ws.Range("J61").Value = 11765508
$ws.Range("L61").Value = 11765508
$ws.Range("N61").Value = -11765932
$ws.Range("H74").Value = 85980056
$ws.Range("I74").Value = 78022630
$ws.Range("J74").Value = 106669350
$ws.Range("K74").Value = 78022630
$ws.Range("L74").Value = 106669350
$ws.Range("M74").Value = -78021756
$ws.Range("N74").Value = -106671098
$ws.Range("H77").Value = 85980056
$ws.Range("I77").Value = 78022630
$ws.Range("J77").Value = 106669350
$ws.Range("K77").Value = 390113150
$ws.Range("L77").Value = 533346750
$ws.Range("M77").Value = -390108782
$ws.Range("N77").Value = -533355486
$ws.Range("H122").Value = 1637.4286
$ws.Range("I122").Value = 820
$ws.Range("J122").Value = 1964.4
$ws.Range("K122").Value = 2460
$ws.Range("L122").Value = 5893.200000000001
$ws.Range("M122").Value = -10
$ws.Range("N122").Value = -10793.2
$ws.Range("H136").Value = 5025761
$ws.Range("J136").Value = 11765508
$ws.Range("L136").Value = 35296524
$ws.Range("N136").Value = -35301624
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 496.48
$ws.Range("I22").Value = 487.91666
$ws.Range("J22").Value = 702
$ws.Range("K22").Value = 487.91666
$ws.Range("L22").Value = 702
$ws.Range("M22").Value = -314.91666
$ws.Range("N22").Value = -1048
$ws.Range("H107").Value = 843.7273
$ws.Range("I107").Value = 593.75
$ws.Range("J107").Value = 986.5714
$ws.Range("K107").Value = 593.75
$ws.Range("L107").Value = 986.5714
$ws.Range("M107").Value = 1326.25
$ws.Range("N107").Value = -4826.5714
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 613.61536
$ws.Range("I5").Value = 247.125
$ws.Range("J5").Value = 1200
$ws.Range("K5").Value = 247.125
$ws.Range("L5").Value = 1200
$ws.Range("M5").Value = -135.125
$ws.Range("N5").Value = -1424
$ws.Range("H58").Value = 3693189.2
$ws.Range("I58").Value = 2105154
$ws.Range("K58").Value = 2105154
$ws.Range("M58").Value = -2104951
$ws.Range("H86").Value = 8651.066000000001
$ws.Range("I86").Value = 13212.546
$ws.Range("J86").Value = 6010.2104
$ws.Range("K86").Value = 13212.546
$ws.Range("L86").Value = 6010.2104
$ws.Range("M86").Value = -12089.546
$ws.Range("N86").Value = -8256.2104
$ws.Range("H89").Value = 8651.066000000001
$ws.Range("I89").Value = 13212.546
$ws.Range("J89").Value = 6010.2104
$ws.Range("K89").Value = 66062.73
$ws.Range("L89").Value = 30051.052
$ws.Range("M89").Value = -60446.73
$ws.Range("N89").Value = -41283.052
$ws.Range("H132").Value = 1854590.2
$ws.Range("I132").Value = 3126964
$ws.Range("K132").Value = 9380892
$ws.Range("M132").Value = -9378362
$ws.Range("H134").Value = 2868113.5
$ws.Range("I134").Value = 14315.75
$ws.Range("J134").Value = 6673177.5
$ws.Range("K134").Value = 42947.25
$ws.Range("L134").Value = 20019532.5
$ws.Range("M134").Value = -40412.25
$ws.Range("N134").Value = -20024602.5
$ws.Range("H136").Value = 3693189.2
$ws.Range("I136").Value = 2105154
$ws.Range("K136").Value = 6315462
$ws.Range("M136").Value = -6312912
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3557463.5
$ws.Range("I5").Value = 3846818.5
$ws.Range("J5").Value = 3031363.2
$ws.Range("K5").Value = 11540455.5
$ws.Range("L5").Value = 9094089.600000001
$ws.Range("M5").Value = -11540343.5
$ws.Range("N5").Value = -9094313.600000001
$ws.Range("H122").Value = 1033.3125
$ws.Range("I122").Value = 313.30768
$ws.Range("J122").Value = 4153.3335
$ws.Range("K122").Value = 2819.76912
$ws.Range("L122").Value = 37380.0015
$ws.Range("M122").Value = -369.7691199999999
$ws.Range("N122").Value = -42280.0015
$ws.Range("H131").Value = 954.1852
$ws.Range("I131").Value = 440.2
$ws.Range("J131").Value = 1071
$ws.Range("K131").Value = 1320.6
$ws.Range("L131").Value = 3213
$ws.Range("M131").Value = 3719.4
$ws.Range("N131").Value = -13293
$ws.Range("H132").Value = 3087.889
$ws.Range("I132").Value = 4502.6665
$ws.Range("J132").Value = 2380.5
$ws.Range("K132").Value = 40523.9985
$ws.Range("L132").Value = 21424.5
$ws.Range("M132").Value = -37993.9985
$ws.Range("N132").Value = -26484.5
$ws.Range("H135").Value = 3557463.5
$ws.Range("I135").Value = 3846818.5
$ws.Range("J135").Value = 3031363.2
$ws.Range("K135").Value = 34621366.5
$ws.Range("L135").Value = 27282268.8
$ws.Range("M135").Value = -34618831.5
$ws.Range("N135").Value = -27287338.8
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11303.48
$ws.Range("I80").Value = 4586.0713
$ws.Range("K80").Value = 4586.0713
$ws.Range("M80").Value = -3588.0713
$ws.Range("H83").Value = 11303.48
$ws.Range("I83").Value = 4586.0713
$ws.Range("K83").Value = 22930.3565
$ws.Range("M83").Value = -17938.3565
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6983.091
$ws.Range("I22").Value = 5750
$ws.Range("J22").Value = 7687.7144
$ws.Range("K22").Value = 5750
$ws.Range("L22").Value = 7687.7144
$ws.Range("M22").Value = -5455
$ws.Range("N22").Value = -8277.714400000001
$ws.Range("H27").Value = 6983.091
$ws.Range("I27").Value = 5750
$ws.Range("J27").Value = 7687.7144
$ws.Range("K27").Value = 5750
$ws.Range("L27").Value = 7687.7144
$ws.Range("M27").Value = -5643
$ws.Range("N27").Value = -7901.7144
$ws.Range("H46").Value = 385.9
$ws.Range("J46").Value = 358.625
$ws.Range("L46").Value = 358.625
$ws.Range("N46").Value = -734.625
$ws.Range("H82").Value = 3998.4546
$ws.Range("J82").Value = 5590.5386
$ws.Range("L82").Value = 5590.5386
$ws.Range("N82").Value = -6312.5386
$ws.Range("H85").Value = 3998.4546
$ws.Range("J85").Value = 5590.5386
$ws.Range("L85").Value = 5590.5386
$ws.Range("N85").Value = -8086.5386
$ws.Range("H132").Value = 1803867
$ws.Range("I132").Value = 2300313
$ws.Range("J132").Value = 4250.125
$ws.Range("K132").Value = 6900939
$ws.Range("L132").Value = 12750.375
$ws.Range("M132").Value = -6898409
$ws.Range("N132").Value = -17810.375
$ws.Range("H136").Value = 13624267
$ws.Range("I136").Value = 28759836
$ws.Range("J136").Value = 2254.6
$ws.Range("K136").Value = 86279508
$ws.Range("L136").Value = 6763.799999999999
$ws.Range("M136").Value = -86276958
$ws.Range("N136").Value = -11863.8
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 12555.4
$ws.Range("I136").Value = 10533.25
$ws.Range("K136").Value = 31599.75
$ws.Range("M136").Value = -29049.75
